$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet to reflect the new export timestamp
$ws.Name = "IClientBalance-20241014-092445-"

# Shift every reference date in column G (rows 2 through 274) from 2024-10-11 to 2024-10-14
for ($r = 2; $r -le 274; $r++) {
    $ws.Cells.Item($r, 7).Value = 45579
}

# Targeted balance corrections (Vl. Projetado / Saldo Previsto / Vl. Total)
$ws.Range("E5").Value = 0
$ws.Range("H5").Value = 0

$ws.Range("D15").Value = -214408.48
$ws.Range("E15").Value = 202727.78
$ws.Range("H15").Value = -11680.7

$ws.Range("E52").Value = 791.79
$ws.Range("H52").Value = 791.79

$ws.Range("E97").Value = 0
$ws.Range("H97").Value = 0

$ws.Range("E101").Value = 0
$ws.Range("H101").Value = 0

$ws.Range("E112").Value = 0.67
$ws.Range("H112").Value = 0.67

$ws.Range("E113").Value = 0.97
$ws.Range("H113").Value = 0.97

$ws.Range("E271").Value = 0
$ws.Range("H271").Value = 0
